# Adding initial omp results from home
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Fill in the new "omp at home" results (columns L, O, R, rows 4-13)
# ---------------------------------------------------------------
$L = @(367572, 367813, 366714, 367798, 365973, 366098, 364059, 366264, 366725, 363865)
$O = @(354232, 353750, 354444, 354779, 353633, 353975, 354619, 354400, 354844, 353943)
$R = @(364832, 362110, 363578, 362158, 364705, 365631, 366666, 364989, 364916, 366954)

for ($i = 0; $i -lt 10; $i++) {
    $row = 4 + $i

    $ws.Range("L$row").Value = $L[$i]
    $ws.Range("L$row").Style = "Normal"

    $ws.Range("O$row").Value = $O[$i]

    $ws.Range("R$row").Value = $R[$i]
}

# ---------------------------------------------------------------
# 2. Re-style the "Parallel For" separator cells on rows 29 & 30
#    (they switch from the plain "Bad" style to the bold "Bad" style,
#     matching the style already used on row 31)
# ---------------------------------------------------------------
$sepCols = @("J", "T", "AD", "AN", "AX", "BG")
foreach ($col in $sepCols) {
    $ws.Range($col + "31").Copy()
    $ws.Range($col + "29").PasteSpecial(-4122)
    $ws.Range($col + "31").Copy()
    $ws.Range($col + "30").PasteSpecial(-4122)
}

# ---------------------------------------------------------------
# 3. Add new rows 32-36 (blank separator rows styled like row 31),
#    with the new commentary text in K33.
# ---------------------------------------------------------------
$blankCols = @("B", "E", "H", "I", "J", "T", "AD", "AN", "AX", "BG")
foreach ($targetRow in 32..36) {
    foreach ($col in $blankCols) {
        $ws.Range($col + "31").Copy()
        $ws.Range($col + ($targetRow)).PasteSpecial(-4122)
    }
    $ws.Rows.Item($targetRow).RowHeight = 23.25
}

$ws.Range("K31").Copy()
$ws.Range("K33").PasteSpecial(-4122)
$ws.Range("K33").Value = "Only Sundaram was sped up at home."

# ---------------------------------------------------------------
# 4. Update the active selection to match the saved workbook state
# ---------------------------------------------------------------
$ws.Range("O32").Select()
